$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Add header row 1: B1:Z1 = 0..24 (sequential column indices)
$cols = @("B","C","D","E","F","G","H","I","J","K","L","M","N","O","P","Q","R","S","T","U","V","W","X","Y","Z")
for ($i = 0; $i -lt $cols.Length; $i++) {
    $ws.Range($cols[$i] + "1").Value = $i
}

# Add column A values for existing rows 2..28 (row index 0..26)
for ($r = 2; $r -le 28; $r++) {
    $ws.Cells.Item($r, 1).Value = $r - 2
}

# Update the selected cell to match the new active selection
$ws.Range("M33").Select()
